$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 169 -- this shifts the existing rows
# 169..184 down to 170..185 (matching the diff's row-by-row shift).
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new record's data.
$ws.Range("A169").Value = 4
$ws.Range("B169").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C169").Value = "Los Lagos"
$ws.Range("D169").Value = 44449
$ws.Range("E169").Value = 10
$ws.Range("F169").Value = 100112023
$ws.Range("G169").Value = "Brócoli"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 1200
$ws.Range("K169").Value = 1300
$ws.Range("L169").Value = 1300
$ws.Range("M169").Value = 1300
$ws.Range("N169").Value = "$/unidad"
$ws.Range("O169").Value = "Región Metropolitana"
$ws.Range("P169").Value = 1300
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = "Hortaliza"

# D column (Fecha) uses a date number format; make sure the freshly
# inserted cell keeps the same formatting as the rest of the Fecha column.
$ws.Range("D169").NumberFormat = $ws.Range("D170").NumberFormat
